# Update countries & provincias Spain
# - Reorders "Rumania" before "Bolivia" and "Islas Malvinas" before "Montserrat"
#   in the shared-string table (achieved here by writing the swapped country
#   names + their corresponding up-to-date stat rows directly).
# - Refreshes the COVID case numbers for a number of countries.
# - Bumps the "Datos actualizados" timestamp from 11:36 to 12:53.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 5 de Octubre de 2020 a las 12:53"

function Set-Row($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# r  A (country)                         B        C     D        E       F  G   H
Set-Row 4   @("Estados Unidos",                 7637066,  154, 4849454, 2572997, 0,  4, 214615)
Set-Row 5   @("India",                          6626291, 4111, 5586703,  936842, 0, 32, 102746)
Set-Row 19  @("Banglades",                       370132, 1442,  283182,   81575, 0, 27,   5375)
Set-Row 31  @("Rumania",                         137491, 1591,  108526,   23917, 0, 45,   5048)
Set-Row 32  @("Bolivia",                         136868,  299,   97547,   31220, 0, 28,   8101)
Set-Row 43  @("Oman",                            101814,  544,   90600,   10229, 0,  8,    985)
Set-Row 44  @("Emiratos Arabes Unidos",           99733,  932,   89410,    9894, 0,  3,    429)
Set-Row 51  @("Bielorrusia",                      80696,  401,   75303,    4531, 0,  5,    862)
Set-Row 62  @("Suiza",                            55932, 1548,   45800,    8054, 0,  1,   2078)
Set-Row 91  @("Madagascar",                       16570,   12,   15601,     736, 0,  1,    233)
Set-Row 93  @("Senegal",                          15122,   28,   12870,    1940, 0,  0,    312)
Set-Row 99  @("Malasia",                          12813,  432,   10340,    2336, 0,  0,    137)
Set-Row 102 @("Consejo Danes para los Refugiados", 10778,   18,   10239,     265, 0,  0,    274)
Set-Row 117 @("Eslovenia",                         6573,   75,    4314,    2103, 0,  1,    156)
Set-Row 127 @("Hong Kong",                         5125,   11,    4864,     156, 0,  0,    105)
Set-Row 144 @("Malta",                             3327,   57,    2770,     518, 0,  0,     39)
Set-Row 181 @("Gibraltar",                          432,    0,     363,      69, 0,  0,      0)
Set-Row 195 @("Liechtenstein",                      127,    1,     116,      10, 0,  0,      1)
Set-Row 215 @("Islas Malvinas",                      13,    0,      13,       0, 0,  0,      0)
Set-Row 216 @("Montserrat",                          13,    0,      12,       0, 0,  0,      1)
